$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "Titile" -> "Title"
$ws.Range("B1").Value = "Title"

# Fix publisher name typo on the "Modern Data Architecture on AWS" row
$ws.Range("E5").Value = "Packt Publishing"

# Add new book rows (10-14)
$ws.Range("A10").Value = "978-1-098-15922-1"
$ws.Range("B10").Value = "Generative AI on AWS"
$ws.Range("C10").Value = "1st"
$ws.Range("D10").Value = 2024
$ws.Range("E10").Value = "O'reilly"
$ws.Range("F10").Value = "Chris Fregly, Antje Barth,Shelbee Eigenbrode"
$ws.Range("G10").Value = "aws,generative ai, large language models,deep learning"

$ws.Range("A11").Value = "978-0387-31073-2"
$ws.Range("B11").Value = "Pattern Recognition and Machine Learning"
$ws.Range("C11").Value = "1st"
$ws.Range("D11").Value = 2006
$ws.Range("E11").Value = "Springer"
$ws.Range("F11").Value = "Christopher M. Bishop"
$ws.Range("G11").Value = "machine learning, ml algorithms, deep learning"

$ws.Range("A12").Value = "978-3-031-45467-7"
$ws.Range("B12").Value = "Deep Learning Foundations and Concepts"
$ws.Range("C12").Value = "1st"
$ws.Range("D12").Value = 2024
$ws.Range("E12").Value = "Springer"
$ws.Range("F12").Value = "Christopher M. Bishop, Hugh Bishop"
$ws.Range("G12").Value = "machine learning, ml algorithms, deep learning"

$ws.Range("A13").Value = "978-1-83763-418-7"
$ws.Range("B13").Value = "15 Math Concepts Every Data Scientist Should Know"
$ws.Range("C13").Value = "1st"
$ws.Range("D13").Value = 2024
$ws.Range("E13").Value = "Packt Publishing"
$ws.Range("F13").Value = "David Hoyle"
$ws.Range("G13").Value = "mathematics, machine learning, deep learning, data science"

$ws.Range("A14").Value = "978-1-83763-111-7"
$ws.Range("B14").Value = "DevSecOps for Azure"
$ws.Range("C14").Value = "1st"
$ws.Range("D14").Value = 2024
$ws.Range("E14").Value = "Packt Publishing"
$ws.Range("F14").Value = "David Okeyode, Joylynn Kirui"
$ws.Range("G14").Value = "devsecops, azure, cloud"

# Update selection to match the saved view state
$ws.Range("B1").Select()
